# Update Sheet4 ("results") with summary description text added to the
# scenario comparison table, per commit "Update table with summary descriptions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

# --- Header row: give the data columns (C1:J1) the wrap-text style that the
# rest of the table already uses, and extend the header across the new
# (currently blank) column K.
$ws.Range("C1:K1").WrapText = $true

# --- Existing rows 2-4: fill in the previously-empty summary cells.
$ws.Range("C2").Value = "More contraction with middle and late century compared to others"
$ws.Range("D2").Value = "Broad but concentrated densities near AK peninsula peaking at mid-century"
$ws.Range("E2").Value = "Reduction over the course of the century"
$ws.Range("F2").Value = "High density near AK peninsula but switches to northern shelf over mid and late century"
$ws.Range("G2").Value = "Western cluster disappears, but pattern in high densities differs for scenario"
$ws.Range("I2").Value = "Reduction over century, high density near Unimak Pass consistent but shrinks"

$ws.Range("C3").Value = "Similar to MIROC"
$ws.Range("D3").Value = "Concentrations near AK peninsula and also north of Pribs"
$ws.Range("E3").Value = "Reduction over the course of the century"
$ws.Range("F3").Value = "Shift from higher densities in south to northern shelf"
$ws.Range("G3").Value = "Similar reduction in density over time for both scenarios"
$ws.Range("I3").Value = "Reduction over century, high density near Unimak Pass consistent but shrinks though less than CESM"

$ws.Range("C4").Value = "Similar to MIROC"
$ws.Range("D4").Value = "Concentrated density at Unimak Pass but some more spread out over the other two hot spots"
$ws.Range("E4").Value = "Reduction over the course of the century"
$ws.Range("F4").Value = "Shift from higher densities in south to northern shelf"
$ws.Range("G4").Value = "Similar to CESM with opposite responses but the western cluster doesn't disappear entirely"
$ws.Range("I4").Value = "Reduction over century, high density near Unimak Pass consistent but shrinks though less than CESM"

# --- Row 5: blank spacer row, but still needs the wrap-text style applied
# across C:K to match the rest of the table.
$ws.Range("C5:K5").WrapText = $true

# --- Row 6: SSP1-2.6 scenario summary.
$ws.Range("A6").Value = "SSP1-2.6"
$ws.Range("B6").Value = "Little change"
$ws.Range("C6").Value = "In some cases, more contraction than high emission scenario"
$ws.Range("D6").Value = "Higher densities toward end of century"
$ws.Range("E6").Value = "Typically smaller area of high density at beginning of century in comparison"
$ws.Range("F6").Value = "Highest densities in northern area show up in late century"
$ws.Range("G6").Value = "Usually a reduction in density over time with western cluster all but disappearing in all cases"
$ws.Range("I6").Value = "Similar amongst scenarios, most differences seen for specific ESMs during mid-century"
$ws.Range("B6:K6").WrapText = $true

# --- Row 7: SSP5-8.5 scenario summary.
$ws.Range("A7").Value = "SSP5-8.5"
$ws.Range("B7").Value = "More contraction of high density areas, particularly in the northern concentration"
$ws.Range("C7").Value = "Sometimes slightly less contraction by end of century"
$ws.Range("D7").Value = "Higher densities at mid-centuries"
$ws.Range("E7").Value = "Wider concentration of high density in comparison, but mid and late century look pretty much the same for both scenarios"
$ws.Range("F7").Value = "High densities start to show up in mid-century"
$ws.Range("G7").Value = "Some increase in density by end of century on mid-shelf"
$ws.Range("I7").Value = "Similar amongst scenarios, most differences seen for specific ESMs during mid-century"
$ws.Range("B7:K7").WrapText = $true

# --- Row 8: blank spacer row, style to match.
$ws.Range("C8:K8").WrapText = $true

# --- Row 9: general notes row.
$ws.Range("A9").Value = "General"
$ws.Range("B9").Value = "Contraction at specific sites (northwestern corner, near Aleutians)"
$ws.Range("C9").Value = "Contraction toward the southeast, further to the inner shelf"
$ws.Range("D9").Value = "Over time there is a spreading of densities across the mid-shelf focused near the AK peninsula but faintly north of the Pribs"
$ws.Range("E9").Value = "Reduction in size of area with high densities of eggs"
$ws.Range("F9").Value = "Increase in area of high densities over course of century but the area near the AK peninsula reduces by quite a bit"
$ws.Range("G9").Value = "Concentrations to the west seem to disappear or reduce over century, high density cluster shows up in some situations, possible problematic isolated cluster below St Matthew"
$ws.Range("H9").Value = "Problems with this one, no patterns showing"
$ws.Range("I9").Value = "Reduction in density across shelf over century but cluster near Unimak Pass stays strong, band of density across the northern shelf seems to narrow over time"
$ws.Range("J9").Value = "Problems with this one, no patterns showing"
$ws.Range("B9:K9").WrapText = $true

# --- Row heights to match Excel's autofit of the newly wrapped text.
$ws.Rows.Item(3).RowHeight = 101.5
$ws.Rows.Item(4).RowHeight = 101.5
$ws.Rows.Item(6).RowHeight = 101.5
$ws.Rows.Item(7).RowHeight = 130.5
$ws.Rows.Item(9).RowHeight = 166

# --- View state: scrolled down with the new last-edited cell selected.
$ws.Range("J6").Select()
$excel.ActiveWindow.ScrollRow = 4
